{"js": "// 1) Remove the stray \"_GoBack\" bookmark that sits just before \"Fin del CU.\"\n//    (it wrapped no text; it's simply deleted from that paragraph).\nconst finResults = context.document.body.search(\"Fin del CU.\", { matchCase: true });\nfinResults.load(\"items\");\nawait context.sync();\n\nif (finResults.items.length > 0) {\n  const finRange = finResults.items[0];\n  const finPara = finRange.paragraphs.getFirst();\n\n  // Rebuild the paragraph exactly as before, just without the bookmark tags.\n  // (deleteBookmark()/collapsed-range delete() aren't reliable for surgically\n  // stripping bookmark markers, so splice in the paragraph's OOXML instead.)\n  const finParaOoxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' +\n    '<w:p w14:paraId=\"3251FDF4\" w14:textId=\"77777777\" w:rsidR=\"00CF34FB\" w:rsidRDefault=\"00CF34FB\" w:rsidP=\"008F3E09\">' +\n    '<w:pPr>' +\n    '<w:pStyle w:val=\"Prrafodelista\"/>' +\n    '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>' +\n    '<w:ind w:hanging=\"198\"/>' +\n    '</w:pPr>' +\n    '<w:r><w:t>Fin del CU.</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  finPara.insertOoxml(finParaOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"CU 019. Consultar Catalogo Proveedor.\" -> \"No aplica\", and the\n//    \"_GoBack\" bookmark now marks this (the last-edited) spot instead.\nconst cuResults = context.document.body.search(\"CU 019. Consultar Catalogo Proveedor.\", { matchCase: true });\ncuResults.load(\"items\");\nawait context.sync();\n\nif (cuResults.items.length > 0) {\n  const cuRange = cuResults.items[0];\n  cuRange.insertText(\"No aplica\", Word.InsertLocation.replace);\n  await context.sync();\n\n  const endOfRun = cuRange.getRange(Word.RangeLocation.end);\n  endOfRun.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the stray \"_GoBack\" bookmark that sits just before \"Fin del CU.\"\n#    (it wrapped no text; it's simply dropped from that paragraph). The\n#    Bookmark/Range Delete() calls don't surgically strip bookmark markers in\n#    this host, so splice the paragraph back in via its OOXML instead.\n$finRange = $d.Content\n$finFound = $finRange.Find.Execute(\"Fin del CU.\")\nif ($finFound) {\n    $finPara = $finRange.Paragraphs(1).Range\n    $finOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData>' + `\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' + `\n        '<w:body>' + `\n        '<w:p w14:paraId=\"3251FDF4\" w14:textId=\"77777777\" w:rsidR=\"00CF34FB\" w:rsidRDefault=\"00CF34FB\" w:rsidP=\"008F3E09\">' + `\n        '<w:pPr>' + `\n        '<w:pStyle w:val=\"Prrafodelista\"/>' + `\n        '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"3\"/></w:numPr>' + `\n        '<w:ind w:hanging=\"198\"/>' + `\n        '</w:pPr>' + `\n        '<w:r><w:t>Fin del CU.</w:t></w:r>' + `\n        '</w:p>' + `\n        '</w:body>' + `\n        '</w:document>' + `\n        '</pkg:xmlData>' + `\n        '</pkg:part>' + `\n        '</pkg:package>'\n    $finPara.InsertXML($finOoxml)\n}\n\n# 2) \"CU 019. Consultar Catalogo Proveedor.\" -> \"No aplica\", and the\n#    \"_GoBack\" bookmark now marks this (the last-edited) spot instead.\n$cuRange = $d.Content\n$cuFound = $cuRange.Find.Execute(\"CU 019. Consultar Catalogo Proveedor.\")\nif ($cuFound) {\n    $cuPara = $cuRange.Paragraphs(1).Range\n    $cuOoxml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData>' + `\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' + `\n        '<w:body>' + `\n        '<w:p w14:paraId=\"43ADB115\" w14:textId=\"77777777\" w:rsidR=\"00654C01\" w:rsidRDefault=\"00CF34FB\" w:rsidP=\"00654C01\">' + `\n        '<w:r><w:t>No aplica</w:t></w:r>' + `\n        '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/>' + `\n        '<w:bookmarkEnd w:id=\"1\"/>' + `\n        '</w:p>' + `\n        '</w:body>' + `\n        '</w:document>' + `\n        '</pkg:xmlData>' + `\n        '</pkg:part>' + `\n        '</pkg:package>'\n    $cuPara.InsertXML($cuOoxml)\n}\n"}
